# 자가채점표(1조).xlsx - update this week's "4인조" (4-member group) self-grading
# sheet: each of the 4 teammates gets credit for the 20-point "게임AI"
# (Game AI) bonus task, participating at 25% each (commit: "1213 B win AI").
#
# Before: K9:K12 (게임AI participation %) were all 0, so L (총합/total,
#         which is (SUM(C:J)*10 + K*20)/25) topped out at 80.
# After:  K9:K12 = 25 each -> L9:L12 (and the K13 column total) recompute to 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4인조")

$ws.Range("K9").Value  = 25
$ws.Range("K10").Value = 25
$ws.Range("K11").Value = 25
$ws.Range("K12").Value = 25
